# Apply the updated crypto price/volume snapshot to Sheet1.
# Column D holds "Price" text that frequently looks numeric (e.g. "1.010",
# "27.901.27"); a leading apostrophe forces Excel to store it as literal
# text (preserving trailing zeros / multi-dot formatting) instead of
# coercing it to a Number, matching the original inlineStr content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.901.27"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'1.888.44"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'337.11"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4692"
$ws.Range("D8").Value = "'0.3973"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'45.95"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.08049"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'1.018"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'22.10"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'1.884.51"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'6.031"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "'7.317"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'89.64"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "'0.06732"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'17.44"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'27.904.94"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'5.522"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'11.06"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "'2.310"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").Value = "'2.114.22"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").Value = "'159.61"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'19.92"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Value = "'2.174"
$ws.Range("E29").Value = "  +3.29%  "
$ws.Range("D30").Value = "'5.533"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'122.30"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'0.9905"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'0.09510"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'3.634"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'5.360"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'1.367"
$ws.Range("E36").Value = "  -5.90%  "
$ws.Range("D37").Value = "'0.06096"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "'0.02255"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'8.403"
$ws.Range("E39").Value = "  +4.28%  "
$ws.Range("D40").Value = "'1.203"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "'1.008"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'0.6023"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'10.46"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5692"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.246"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'12.35"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").Value = "'0.06799"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").Value = "'112.93"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'3.043"
$ws.Range("E51").Value = "  -10.39%  "
